# Automatische test-sync: 2025-06-18 08:00:10
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append the new log entry as row 67
$newRow = 67
$logs.Cells.Item($newRow, 1).Value = "Vragen over samenwerking"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 6).Value = "2025-06-17 23:04:10"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the conditional formatting ranges to cover the newly added row
$dFormats = $logs.Range("D2:D66").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D67"))
}

$gFormats = $logs.Range("G2:G66").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G67"))
}

# Update the Dashboard count for "Overig" category (row 3, column B)
$dashboard.Cells.Item(3, 2).Value = 20
